# Refresh cryptos list: update prices / 1h volume deltas, and restore the
# correct Litecoin/SuiNetwork and NEARProtocol/RenderToken row order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.595.14"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.386.38"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.63"
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.388.86"
$ws.Range("E9").Value = "  +1.99%  "
$ws.Range("E10").Value = "  +5.42%  "
$ws.Range("E11").Value = "  +1.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.37"
$ws.Range("E12").Value = "  +3.09%  "
$ws.Range("E13").Value = "  +4.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.56"
$ws.Range("E14").Value = "  +3.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000169"
$ws.Range("E15").Value = "  +5.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.507.73"
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.05"
$ws.Range("E17").Value = "  +5.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "321.97"
$ws.Range("E18").Value = "  +2.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.15"
$ws.Range("E19").Value = "  +2.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.79"
$ws.Range("E20").Value = "  +4.37%  "
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("B22").Value = "SuiNetwork"
$ws.Range("C22").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.75"
$ws.Range("E22").Value = "  -5.90%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.49"
$ws.Range("E23").Value = "  +2.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.96"
$ws.Range("E24").Value = "  +11.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.27"
$ws.Range("E25").Value = "  +4.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "524.36"
$ws.Range("E26").Value = "  +3.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0912"
$ws.Range("E27").Value = "  +2.56%  "
$ws.Range("E28").Value = "  +5.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.40"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("E30").Value = "  +2.86%  "
$ws.Range("E31").Value = "  +2.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.60"
$ws.Range("E33").Value = "  +6.95%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.75"
$ws.Range("E34").Value = "  +5.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.91"
$ws.Range("E35").Value = "  +7.78%  "
$ws.Range("E36").Value = "  +2.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "147.31"
$ws.Range("E38").Value = "  +6.08%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "41.42"
$ws.Range("E40").Value = "  +3.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "148.96"
$ws.Range("E41").Value = "  +8.78%  "
$ws.Range("E42").Value = "  +5.60%  "
$ws.Range("E43").Value = "  +2.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0527"
$ws.Range("E44").Value = "  +3.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.88"
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.585"
$ws.Range("E46").Value = "  +3.54%  "
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("E48").Value = "  +1.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.41"
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.84"
$ws.Range("E50").Value = "  +2.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.02"
$ws.Range("E51").Value = "  +5.05%  "
